# "moving to new server" — update the label and figure on Лист1, and
# move the active selection to B1 (where the figure lives).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

$ws.Range("A1").Value = "Центр"
$ws.Range("B1").Value = 33

$ws.Range("B1").Select()
